$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.424.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "'1.944.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("D5").Value = "'242.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'57.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").Value = "'0.360"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("D10").Value = "'0.0852"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").Value = "'0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "'2.228.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "'0.811"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.49%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'21.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'13.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "'5.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.94%  "
$ws.Range("D17").Value = "'1.951.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "'36.373.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "'69.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "'0.0₃0863"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.79%  "
$ws.Range("D21").Value = "'228.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").Value = "'4.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.48%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E24").Value = "  -5.53%  "
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").Value = "'9.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.37%  "
$ws.Range("D27").Value = "'161.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.131"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  -6.24%  "
$ws.Range("D32").Value = "'4.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.51%  "
$ws.Range("D33").Value = "'0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("D34").Value = "'4.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("D35").Value = "'6.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.04%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").Value = "'2.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "'3.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.58%  "
$ws.Range("D40").Value = "'0.0984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'0.0210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").Value = "'15.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "'1.340.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  -5.44%  "
$ws.Range("D47").Value = "'86.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("D48").Value = "'7.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'2.119.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "'42.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.75%  "
